$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"
$ws.Range("E1").Value = "CITY"
$ws.Range("F1").Value = "STATE"
$ws.Range("G1").Value = "ZIP CODE"
$ws.Range("H1").Value = "PHONE #"
$ws.Range("I1").Value = "PROJECT COORDINATOR"
$ws.Range("J1").Value = "CERTIFICATION REQUIRED (Yes or No)"
$ws.Range("K1").Value = "CERTIFICATION COORDINATOR"

$ws.Range("B2").Value = "YYY"
$ws.Range("C2").Value = "123ABX007"
$ws.Range("D2").Value = "Karapakkam"
$ws.Range("E2").Value = "Chennai"
$ws.Range("F2").Value = "Tamil Nadu"
$ws.Range("G2").Value = "'600117"
$ws.Range("H2").Value = "'9911991100"
$ws.Range("I2").Value = "Sam"
$ws.Range("J2").Value = "Yes"
$ws.Range("K2").Value = "'"
